$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A28").Font.Name = "Calibri"
